$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -6.049999999999998
$ws.Range("D4").Value = -7.538899999999999
$ws.Range("D7").Value = -8.227099999999993
$ws.Range("D8").Value = -8.562499999999996
$ws.Range("B11").Value = 5.399700000000001
$ws.Range("B12").Value = 5.325799999999999
$ws.Range("D12").Value = -8.109399999999999
$ws.Range("D14").Value = -8.593999999999999
$ws.Range("B15").Value = 5.020000000000002
$ws.Range("D22").Value = -7.466899999999996
